# Restructure "饮料" (beverage) sheet:
#  1. For every year block (rows laid out as quarter A/B/C/D), the 2nd and
#     3rd quarter rows (the "B" and "C" rows) had their data rows swapped
#     relative to their quarter label -- i.e. row that used to hold the "B"
#     quarter's figures now holds the "C" quarter's figures and vice versa.
#  2. Columns F ("饮料产销率") and G ("饮料销售量") are removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row pairs (B-quarter row, C-quarter row) to swap, one pair per year
# (years 2000-2019, 4 rows per year starting at row 2).
$pairs = @(
    @(3, 4),
    @(7, 8),
    @(11, 12),
    @(15, 16),
    @(19, 20),
    @(23, 24),
    @(27, 28),
    @(31, 32),
    @(35, 36),
    @(39, 40),
    @(43, 44),
    @(47, 48),
    @(51, 52),
    @(55, 56),
    @(59, 60),
    @(63, 64),
    @(67, 68),
    @(71, 72),
    @(75, 76),
    @(79, 80)
)

foreach ($p in $pairs) {
    $r1 = $p[0]
    $r2 = $p[1]
    $rng1 = $ws.Range("A" + $r1 + ":E" + $r1)
    $rng2 = $ws.Range("A" + $r2 + ":E" + $r2)
    $v1 = $rng1.Value()
    $v2 = $rng2.Value()
    $rng1.Value = $v2
    $rng2.Value = $v1
}

# Drop the now-unneeded "饮料产销率" (F) and "饮料销售量" (G) columns.
$ws.Columns("F:G").Delete()

"ok"
